$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.068178653717041
$ws.Range("B1").Value = 2.427021265029907
$ws.Range("C1").Value = 5.110594272613525
$ws.Range("D1").Value = 2.286472320556641
$ws.Range("E1").Value = 1.304511547088623
